# Apply commit "removed some of the variables and renamed them" to the
# "env" worksheet (and adjust the active-sheet/selection bookkeeping that
# Excel persists alongside it).

$wb  = $excel.ActiveWorkbook
$env = $wb.Worksheets.Item("env")
$spe = $wb.Worksheets.Item("spe")

# --- Remove the pH_KCl column (14th / "N") -----------------------------
# Shifts Ca..Total_C (old O:Y) one column to the left (new N:X).
$env.Range("N1:N19").EntireColumn.Delete()

# --- Remove the trailing k_decomp / S_decomp columns (now Y:Z) --------
$env.Range("Y1:Z19").EntireColumn.Delete()

# --- Rename the remaining header labels --------------------------------
# New shared-string insertion order must be elevation, slope, heatload,
# folded_aspect, soil_depth, so touch the columns in that order.
$env.Cells.Item(1, 2).Value  = "elevation"       # B1 Elevation      -> elevation
$env.Cells.Item(1, 4).Value  = "slope"           # D1 Slope          -> slope
$env.Cells.Item(1, 5).Value  = "heatload"         # E1 Heatload       -> heatload
$env.Cells.Item(1, 3).Value  = "folded_aspect"   # C1 Aspect_folded  -> folded_aspect
$env.Cells.Item(1, 12).Value = "soil_depth"      # L1 Soil_depth     -> soil_depth

# --- Restore view bookkeeping: "env" had focus/column-N selected; the
#     edit leaves "spe" as the active tab, with "env" remembering a
#     full-column selection on N (where pH_KCl used to live). -----------
$env.Range("N1:N1048576").Select()
$spe.Activate()
